$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.7122
$ws.Range("B10").Value = 8.623400000000007
$ws.Range("B12").Value = 5.698000000000003
$ws.Range("B18").Value = 5.008100000000004
